$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.153770666666667
$ws.Range("H2").Value = 3.461312
$ws.Range("I2").Value = 0.1294297218267158
$ws.Range("J2").Value = 0.1294297218267158
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.248835333333334
$ws.Range("N2").Value = 6.746506
$ws.Range("O2").Value = 0.03590294220158827
$ws.Range("P2").Value = 0.03590294220158827
$ws.Range("Q2").Value = 2.594640241763555
$ws.Range("R2").Value = 23.351762175872
$ws.Range("S2").Value = 0.004646907821912227
$ws.Range("T2").Value = 0.004646907821912225
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.153770666666667
$ws.Range("H3").Value = 3.461312
$ws.Range("I3").Value = 0.1294297218267158
$ws.Range("J3").Value = 0.1294297218267158
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 44.29005966666667
$ws.Range("N3").Value = 132.870179
$ws.Range("O3").Value = 0.7070964373190639
$ws.Range("P3").Value = 0.7070964373190639
$ws.Range("Q3").Value = 51.10057166831644
$ws.Range("R3").Value = 459.905145014848
$ws.Range("S3").Value = 0.09151929518686824
$ws.Range("T3").Value = 0.09151929518686823
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.153770666666667
$ws.Range("H4").Value = 3.461312
$ws.Range("I4").Value = 0.1294297218267158
$ws.Range("J4").Value = 0.1294297218267158
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.09762433333333
$ws.Range("N4").Value = 48.292873
$ws.Range("O4").Value = 0.2570006204793478
$ws.Range("P4").Value = 0.2570006204793479
$ws.Range("Q4").Value = 18.57296675881955
$ws.Range("R4").Value = 167.156700829376
$ws.Range("S4").Value = 0.03326351881793535
$ws.Range("T4").Value = 0.03326351881793535
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.560947
$ws.Range("H5").Value = 1.682841
$ws.Range("I5").Value = 0.06292690243138796
$ws.Range("J5").Value = 0.06292690243138795
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.248835333333334
$ws.Range("N5").Value = 6.746506
$ws.Range("O5").Value = 0.03590294220158827
$ws.Range("P5").Value = 0.03590294220158827
$ws.Range("Q5").Value = 1.261477433727333
$ws.Range("R5").Value = 11.353296903546
$ws.Range("S5").Value = 0.002259260940919106
$ws.Range("T5").Value = 0.002259260940919106
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.560947
$ws.Range("H6").Value = 1.682841
$ws.Range("I6").Value = 0.06292690243138796
$ws.Range("J6").Value = 0.06292690243138795
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.29005966666667
$ws.Range("N6").Value = 132.870179
$ws.Range("O6").Value = 0.7070964373190639
$ws.Range("P6").Value = 0.7070964373190639
$ws.Range("Q6").Value = 24.84437609983767
$ws.Range("R6").Value = 223.599384898539
$ws.Range("S6").Value = 0.04449538852075877
$ws.Range("T6").Value = 0.04449538852075876
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.560947
$ws.Range("H7").Value = 1.682841
$ws.Range("I7").Value = 0.06292690243138796
$ws.Range("J7").Value = 0.06292690243138795
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.09762433333333
$ws.Range("N7").Value = 48.292873
$ws.Range("O7").Value = 0.2570006204793478
$ws.Range("P7").Value = 0.2570006204793479
$ws.Range("Q7").Value = 9.029914076910332
$ws.Range("R7").Value = 81.269226692193
$ws.Range("S7").Value = 0.01617225296971008
$ws.Range("T7").Value = 0.01617225296971008
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.837574333333334
$ws.Range("H8").Value = 11.512723
$ws.Range("I8").Value = 0.4304981854736105
$ws.Range("J8").Value = 0.4304981854736104
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.248835333333334
$ws.Range("N8").Value = 6.746506
$ws.Range("O8").Value = 0.03590294220158827
$ws.Range("P8").Value = 0.03590294220158827
$ws.Range("Q8").Value = 8.630072755093114
$ws.Range("R8").Value = 77.670654795838
$ws.Range("S8").Value = 0.01545615147094767
$ws.Range("T8").Value = 0.01545615147094766
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.837574333333334
$ws.Range("H9").Value = 11.512723
$ws.Range("I9").Value = 0.4304981854736105
$ws.Range("J9").Value = 0.4304981854736104
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 44.29005966666667
$ws.Range("N9").Value = 132.870179
$ws.Range("O9").Value = 0.7070964373190639
$ws.Range("P9").Value = 0.7070964373190639
$ws.Range("Q9").Value = 169.9663961986019
$ws.Range("R9").Value = 1529.697565787417
$ws.Range("S9").Value = 0.3044037332207116
$ws.Range("T9").Value = 0.3044037332207115
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.837574333333334
$ws.Range("H10").Value = 11.512723
$ws.Range("I10").Value = 0.4304981854736105
$ws.Range("J10").Value = 0.4304981854736104
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.09762433333333
$ws.Range("N10").Value = 48.292873
$ws.Range("O10").Value = 0.2570006204793478
$ws.Range("P10").Value = 0.2570006204793479
$ws.Range("Q10").Value = 61.77582996924212
$ws.Range("R10").Value = 555.9824697231791
$ws.Range("S10").Value = 0.1106383007819512
$ws.Range("T10").Value = 0.1106383007819512
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.179188666666667
$ws.Range("H11").Value = 9.537566
$ws.Range("I11").Value = 0.3566406363494371
$ws.Range("J11").Value = 0.356640636349437
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.248835333333334
$ws.Range("N11").Value = 6.746506
$ws.Range("O11").Value = 0.03590294220158827
$ws.Range("P11").Value = 0.03590294220158827
$ws.Range("Q11").Value = 7.149471804932889
$ws.Range("R11").Value = 64.34524624439599
$ws.Range("S11").Value = 0.0128044481535915
$ws.Range("T11").Value = 0.0128044481535915
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.179188666666667
$ws.Range("H12").Value = 9.537566
$ws.Range("I12").Value = 0.3566406363494371
$ws.Range("J12").Value = 0.356640636349437
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.29005966666667
$ws.Range("N12").Value = 132.870179
$ws.Range("O12").Value = 0.7070964373190639
$ws.Range("P12").Value = 0.7070964373190639
$ws.Range("Q12").Value = 140.8064557382571
$ws.Range("R12").Value = 1267.258101644314
$ws.Range("S12").Value = 0.2521793233658908
$ws.Range("T12").Value = 0.2521793233658908
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.179188666666667
$ws.Range("H13").Value = 9.537566
$ws.Range("I13").Value = 0.3566406363494371
$ws.Range("J13").Value = 0.356640636349437
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 16.09762433333333
$ws.Range("N13").Value = 48.292873
$ws.Range("O13").Value = 0.2570006204793478
$ws.Range("P13").Value = 0.2570006204793479
$ws.Range("Q13").Value = 51.17738484079089
$ws.Range("R13").Value = 460.596463567118
$ws.Range("S13").Value = 0.09165686482995478
$ws.Range("T13").Value = 0.09165686482995478
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.182783
$ws.Range("H14").Value = 0.548349
$ws.Range("I14").Value = 0.02050455391884863
$ws.Range("J14").Value = 0.02050455391884863
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.248835333333334
$ws.Range("N14").Value = 6.746506
$ws.Range("O14").Value = 0.03590294220158827
$ws.Range("P14").Value = 0.03590294220158827
$ws.Range("Q14").Value = 0.4110488687326667
$ws.Range("R14").Value = 3.699439818594
$ws.Range("S14").Value = 0.0007361738142177729
$ws.Range("T14").Value = 0.0007361738142177728
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.182783
$ws.Range("H15").Value = 0.548349
$ws.Range("I15").Value = 0.02050455391884863
$ws.Range("J15").Value = 0.02050455391884863
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 44.29005966666667
$ws.Range("N15").Value = 132.870179
$ws.Range("O15").Value = 0.7070964373190639
$ws.Range("P15").Value = 0.7070964373190639
$ws.Range("Q15").Value = 8.095469976052334
$ws.Range("R15").Value = 72.859229784471
$ws.Range("S15").Value = 0.01449869702483452
$ws.Range("T15").Value = 0.01449869702483452
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.182783
$ws.Range("H16").Value = 0.548349
$ws.Range("I16").Value = 0.02050455391884863
$ws.Range("J16").Value = 0.02050455391884863
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 16.09762433333333
$ws.Range("N16").Value = 48.292873
$ws.Range("O16").Value = 0.2570006204793478
$ws.Range("P16").Value = 0.2570006204793479
$ws.Range("Q16").Value = 2.942372068519667
$ws.Range("R16").Value = 26.481348616677
$ws.Range("S16").Value = 0.005269683079796341
$ws.Range("T16").Value = 0.005269683079796342
